{"js": "// Applies the style changes described by the diff:\n//  1. Add a new paragraph style \"Abstract Title\" (styleId AbstractTitle)\n//  2. Change the \"Abstract\" style's space-before from 300 (15pt) to 100 (5pt)\n//  3. Give the \"ImportTok\" character style a green, bold color\n//  4. Give the \"BuiltInTok\" character style a green color\n\nconst styles = context.document.getStyles();\n\n// --- 1. New \"Abstract Title\" paragraph style -------------------------------\n// Create it first and sync so the new style actually exists in the\n// document before we try to configure it (a freshly-minted style's proxy\n// cannot take further property writes until it is re-fetched by name).\ncontext.document.addStyle(\"Abstract Title\", Word.StyleType.paragraph);\nawait context.sync();\n\nconst abstractTitle = styles.getByName(\"Abstract Title\");\nabstractTitle.baseStyle = \"Normal\";\nabstractTitle.nextParagraphStyle = \"Abstract\";\nabstractTitle.quickStyle = true;\n\nabstractTitle.paragraphFormat.keepWithNext = true;\nabstractTitle.paragraphFormat.keepTogether = true;\nabstractTitle.paragraphFormat.alignment = Word.Alignment.centered;\nabstractTitle.paragraphFormat.spaceBefore = 15;\nabstractTitle.paragraphFormat.spaceAfter = 0;\n\nabstractTitle.font.size = 10;\nabstractTitle.font.sizeBidirectional = 10;\nabstractTitle.font.bold = true;\nabstractTitle.font.color = \"#345A8A\";\n\n// --- 2. \"Abstract\" style: spacing before 300 -> 100 (twips), i.e. 15pt -> 5pt\nconst abstractStyle = styles.getByName(\"Abstract\");\nabstractStyle.paragraphFormat.spaceBefore = 5;\n\n// --- 3. \"ImportTok\" character style: add bold + green color ----------------\nconst importTok = styles.getByName(\"ImportTok\");\nimportTok.font.color = \"#008000\";\nimportTok.font.bold = true;\n\n// --- 4. \"BuiltInTok\" character style: add green color -----------------------\nconst builtInTok = styles.getByName(\"BuiltInTok\");\nbuiltInTok.font.color = \"#008000\";\n\nawait context.sync();\n", "ps1": "# Word COM interop script applying the style changes described by the diff:\n#  1. Add a new paragraph style \"Abstract Title\" (styleId AbstractTitle)\n#  2. Change the \"Abstract\" style's space-before from 300 (15pt) to 100 (5pt)\n#  3. Give the \"ImportTok\" character style a green, bold color\n#  4. Give the \"BuiltInTok\" character style a green color\n\n$d = $word.ActiveDocument\n\n# --- 1. New \"Abstract Title\" paragraph style -------------------------------\n$abstractTitle = $d.Styles.Add(\"Abstract Title\", 1)\n$abstractTitle.BaseStyle = $d.Styles(\"Normal\")\n$abstractTitle.NextParagraphStyle = $d.Styles(\"Abstract\")\n$abstractTitle.QuickStyle = $true\n\n$abstractTitle.ParagraphFormat.KeepWithNext = $true\n$abstractTitle.ParagraphFormat.KeepTogether = $true\n$abstractTitle.ParagraphFormat.Alignment = 1\n$abstractTitle.ParagraphFormat.SpaceBefore = 15\n$abstractTitle.ParagraphFormat.SpaceAfter = 0\n\n$abstractTitle.Font.Size = 10\n$abstractTitle.Font.SizeBi = 10\n$abstractTitle.Font.Bold = $true\n$abstractTitle.Font.Color = 9067060\n\n# --- 2. \"Abstract\" style: spacing before 300 -> 100 ------------------------\n$abstract = $d.Styles(\"Abstract\")\n$abstract.ParagraphFormat.SpaceBefore = 5\n\n# --- 3. \"ImportTok\" character style: add bold + green color ----------------\n$importTok = $d.Styles(\"ImportTok\")\n$importTok.Font.Color = 32768\n$importTok.Font.Bold = $true\n\n# --- 4. \"BuiltInTok\" character style: add green color -----------------------\n$builtInTok = $d.Styles(\"BuiltInTok\")\n$builtInTok.Font.Color = 32768\n"}
